# Add two new columns, I ("I0") and J ("IF"), to the existing data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, thin border, centered alignment)
# from an existing header cell (H1) onto the two new header cells so the
# new headers pick up the same cell style used by the rest of row 1,
# instead of creating a brand-new style.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: column I mirrors a fixed "innings pitched" style value (1,
# except for the very first data row which is 7), column J mirrors the
# existing column H value for that row.
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 5

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 5

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 2

$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 3

$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 3
